$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1 / Worksheets index 1): update column F values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 712
$wsExhibit.Range("F4").Value = 240
$wsExhibit.Range("F5").Value = 2502
$wsExhibit.Range("F6").Value = 54
$wsExhibit.Range("F7").Value = 3591
$wsExhibit.Range("F9").Value = 910

# Sheet "全部类型" (sheet4 / Worksheets index 4): update column F values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 712
$wsAll.Range("F5").Value = 240
$wsAll.Range("F6").Value = 2502
$wsAll.Range("F7").Value = 54
$wsAll.Range("F8").Value = 3591
$wsAll.Range("F10").Value = 910
